# GridMaze maze_hub_1_1 BOM update
#
# - A new IR illumination connector (J1, Molex Micro-Fit connector,
#   part 43045-0210) is added to the BOM as a new row, pushing every
#   following row down by one.
# - The existing "link to pyboard store" hyperlink (originally on F36)
#   naturally ends up on F37 once the new row is inserted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 9 (SV1 / Right angle receptacle),
# shifting all subsequent rows (and their hyperlink) down by one.
$ws.Rows(9).Insert()

# Populate the newly inserted row 9 with the Micro-Fit connector part.
$ws.Range("A9").Value = "J1"
$ws.Range("B9").Value = "-"
$ws.Range("C9").Value = "-"
$ws.Range("D9").Value = "-"
$ws.Range("E9").Value = "Micro-Fit connector"
$ws.Range("F9").Value = "Molex"
$ws.Range("G9").Value = "43045-0210"
$ws.Range("H9").Value = 3103032

# The hyperlink that used to live on F36 needs to move down to F37 along
# with the row of data it decorates (this engine's Rows.Insert does not
# automatically re-anchor existing Hyperlink objects), so recreate it and
# restore the cell's original "Hyperlink" look (Add() nudges the cell onto
# a freshly minted style record otherwise).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F37"), "https://micropython.org/store")
$ws.Range("F37").Style = "Hyperlink"

# Match the author's final cursor/selection position recorded in the file.
$ws.Range("N21").Select()
